$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,92,95,96,98,99,100,103,105,106,111,113,114,115,123,125,126,127,129,130,134,136,139,141,142,144,152,165,174,187,197,199,290)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = 1.00
}
